# Sprint 2 daily Burndown Chart update.
# Sheet1 holds the burndown data (Calendar Days / Planned / Actual);
# the sprint has progressed one more day, so the "Actual" series picks up
# its next data point (C10 = 18, matching the existing idx=7 run-rate).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Record today's "Actual" burndown value for day 9 (row 10).
$ws1.Range("C10").Value = 18

# Nudge the charts to pick up the newly-entered data point. Both the
# Sheet1 chart and the Sheet2 chart plot Sheet1!$C$2:$C$23, so re-assert
# each series' formula (same reference, now covering the fresh cell) to
# prompt a cache refresh against the updated range.
$actualFormula = "=SERIES(Sheet1!`$C`$1,Sheet1!`$A`$2:`$A`$23,Sheet1!`$C`$2:`$C`$23,2)"
foreach ($wsName in @("Sheet1", "Sheet2")) {
    $sheet = $wb.Worksheets.Item($wsName)
    $chartObjects = $sheet.ChartObjects()
    for ($i = 1; $i -le $chartObjects.Count; $i++) {
        $chart = $chartObjects.Item($i).Chart
        $series = $chart.SeriesCollection()
        if ($series.Count -ge 2) {
            $series.Item(2).Formula = $actualFormula
        }
    }
}

# Move the active selection to reflect where work continued (D12), as
# recorded by the workbook's last saved view state.
$ws1.Range("D12").Select()
